$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = [double]"0.09809395674285802"
$ws.Range("B3").Value = [double]"0.002096585386780465"
$ws.Range("C3").Value = [double]"0.0006433701386290832"
$ws.Range("D3").Value = [double]"2.479265742663204"
$ws.Range("E3").Value = [double]"0.05622613508614493"
$ws.Range("F3").Value = [double]"0.0008355991004943173"
$ws.Range("G3").Value = [double]"0.003357571673066612"
$ws.Range("H3").Value = [double]"0.1001905421296385"
$ws.Range("B4").Value = [double]"0.003736332474924525"
$ws.Range("C4").Value = [double]"0.001444033473920208"
$ws.Range("D4").Value = [double]"1.668123128662512"
$ws.Range("E4").Value = [double]"0.0665154353468555"
$ws.Range("F4").Value = [double]"0.0009060699831735727"
$ws.Range("G4").Value = [double]"0.006566594966675479"
$ws.Range("H4").Value = [double]"0.1018302892177825"
$ws.Range("B5").Value = [double]"0.01125870718517311"
$ws.Range("C5").Value = [double]"0.003344418885552435"
$ws.Range("D5").Value = [double]"3.517001549289796"
$ws.Range("E5").Value = [double]"0.09458201444407748"
$ws.Range("F5").Value = [double]"0.004703745370455524"
$ws.Range("G5").Value = [double]"0.0178136689998907"
$ws.Range("H5").Value = [double]"0.1093526639280311"
$ws.Range("B6").Value = [double]"0.003594862653061924"
$ws.Range("C6").Value = [double]"0.004248607003150161"
$ws.Range("D6").Value = [double]"-1.181590339607812"
$ws.Range("E6").Value = [double]"0.07375754654773362"
$ws.Range("F6").Value = [double]"-0.004732281695297018"
$ws.Range("G6").Value = [double]"0.01192200700142087"
$ws.Range("H6").Value = [double]"0.1016888193959199"
$ws.Range("B7").Value = [double]"0.00876852859107111"
$ws.Range("C7").Value = [double]"0.004894371280356733"
$ws.Range("D7").Value = [double]"0.8227831804600985"
$ws.Range("E7").Value = [double]"0.110605291553067"
$ws.Range("F7").Value = [double]"-0.0008242938031011805"
$ws.Range("G7").Value = [double]"0.0183613509852434"
$ws.Range("H7").Value = [double]"0.1068624853339291"
$ws.Range("B8").Value = [double]"0.009869480978140737"
$ws.Range("C8").Value = [double]"0.004837372719276076"
$ws.Range("D8").Value = [double]"1.692786369697473"
$ws.Range("E8").Value = [double]"0.1208621064108152"
$ws.Range("F8").Value = [double]"0.0003883746541285667"
$ws.Range("G8").Value = [double]"0.0193505873021529"
$ws.Range("H8").Value = [double]"0.1079634377209988"
$ws.Range("B9").Value = [double]"0.005933409132903259"
$ws.Range("C9").Value = [double]"0.007215175663447918"
$ws.Range("D9").Value = [double]"2.992229391224369"
$ws.Range("E9").Value = [double]"0.1347819007092562"
$ws.Range("F9").Value = [double]"-0.008208119123971964"
$ws.Range("G9").Value = [double]"0.02007493738977849"
$ws.Range("H9").Value = [double]"0.1040273658757613"
$ws.Range("B10").Value = [double]"-0.09809395674285802"
$ws.Range("C10").Value = [double]"0.00049680731378803"
$ws.Range("D10").Value = [double]"-218.1829518639617"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.09906768426841009"
$ws.Range("G10").Value = [double]"-0.09712022921730593"
$ws.Range("B11").Value = [double]"-0.04484420884934058"
$ws.Range("C11").Value = [double]"0.0008301602795880798"
$ws.Range("D11").Value = [double]"-87.99768868286277"
$ws.Range("E11").Value = [double]"1.904552713363108e-10"
$ws.Range("F11").Value = [double]"-0.04647133533187821"
$ws.Range("G11").Value = [double]"-0.04321708236680295"
$ws.Range("H11").Value = [double]"0.05324974789351744"
$ws.Range("B12").Value = [double]"-0.03397284502487496"
$ws.Range("C12").Value = [double]"0.0006861649245026288"
$ws.Range("D12").Value = [double]"-79.10198495644318"
$ws.Range("E12").Value = [double]"1.946797866040948e-35"
$ws.Range("F12").Value = [double]"-0.03531773000210383"
$ws.Range("G12").Value = [double]"-0.0326279600476461"
$ws.Range("H12").Value = [double]"0.06412111171798306"
$ws.Range("B13").Value = [double]"-0.0252342452223311"
$ws.Range("C13").Value = [double]"0.0006224580240111563"
$ws.Range("D13").Value = [double]"-62.50428019951165"
$ws.Range("E13").Value = [double]"1.144723384872349e-29"
$ws.Range("F13").Value = [double]"-0.02645425907822476"
$ws.Range("G13").Value = [double]"-0.02401423136643743"
$ws.Range("H13").Value = [double]"0.07285971152052692"
$ws.Range("B14").Value = [double]"-0.02005534765890222"
$ws.Range("C14").Value = [double]"0.0006070190706090034"
$ws.Range("D14").Value = [double]"-53.9503442550001"
$ws.Range("E14").Value = [double]"1.840139101163723e-78"
$ws.Range("F14").Value = [double]"-0.02124510068313133"
$ws.Range("G14").Value = [double]"-0.01886559463467309"
$ws.Range("H14").Value = [double]"0.07803860908395581"
$ws.Range("B15").Value = [double]"-0.01569468220469097"
$ws.Range("C15").Value = [double]"0.0005962386375099735"
$ws.Range("D15").Value = [double]"-44.27123054175959"
$ws.Range("E15").Value = [double]"1.317155125655248e-07"
$ws.Range("F15").Value = [double]"-0.01686330590226877"
$ws.Range("G15").Value = [double]"-0.01452605850711317"
$ws.Range("H15").Value = [double]"0.08239927453816706"
$ws.Range("B16").Value = [double]"-0.01460480334155488"
$ws.Range("C16").Value = [double]"0.0005951079454677329"
$ws.Range("D16").Value = [double]"-41.05663540327294"
$ws.Range("E16").Value = [double]"0.0001951048819665571"
$ws.Range("F16").Value = [double]"-0.01577121097624246"
$ws.Range("G16").Value = [double]"-0.0134383957068673"
$ws.Range("H16").Value = [double]"0.08348915340130314"
$ws.Range("B17").Value = [double]"-0.01479829100638238"
$ws.Range("C17").Value = [double]"0.0005985832914485713"
$ws.Range("D17").Value = [double]"-42.76001157318437"
$ws.Range("E17").Value = [double]"3.141664115688989e-18"
$ws.Range("F17").Value = [double]"-0.01597150956612792"
$ws.Range("G17").Value = [double]"-0.01362507244663685"
$ws.Range("H17").Value = [double]"0.08329566573647564"
$ws.Range("B18").Value = [double]"-0.01180063871612091"
$ws.Range("C18").Value = [double]"0.0006029509276619404"
$ws.Range("D18").Value = [double]"-34.50175945168449"
$ws.Range("E18").Value = [double]"1.681921628520686e-16"
$ws.Range("F18").Value = [double]"-0.01298241806673647"
$ws.Range("G18").Value = [double]"-0.01061885936550535"
$ws.Range("H18").Value = [double]"0.08629331802673711"
$ws.Range("B19").Value = [double]"-0.009701550481288632"
$ws.Range("C19").Value = [double]"0.0006019723660431902"
$ws.Range("D19").Value = [double]"-27.72222570072581"
$ws.Range("E19").Value = [double]"0.04632732530615839"
$ws.Range("F19").Value = [double]"-0.01088141204016686"
$ws.Range("G19").Value = [double]"-0.008521688922410403"
$ws.Range("H19").Value = [double]"0.08839240626156938"
$ws.Range("B20").Value = [double]"-0.008032287878410796"
$ws.Range("C20").Value = [double]"0.0006145504994456182"
$ws.Range("D20").Value = [double]"-22.61872265071058"
$ws.Range("E20").Value = [double]"6.71871974089773e-06"
$ws.Range("F20").Value = [double]"-0.009236802584384711"
$ws.Range("G20").Value = [double]"-0.006827773172436889"
$ws.Range("H20").Value = [double]"0.09006166886444722"
$ws.Range("B21").Value = [double]"-0.005584451525468486"
$ws.Range("C21").Value = [double]"0.0006400550324503538"
$ws.Range("D21").Value = [double]"-14.96043862262047"
$ws.Range("E21").Value = [double]"0.02277132466252779"
$ws.Range("F21").Value = [double]"-0.006838956609867467"
$ws.Range("G21").Value = [double]"-0.004329946441069502"
$ws.Range("H21").Value = [double]"0.09250950521738953"
$ws.Range("B22").Value = [double]"-0.002467335563972208"
$ws.Range("C22").Value = [double]"0.0007287223570986099"
$ws.Range("D22").Value = [double]"-8.604089954976807"
$ws.Range("E22").Value = [double]"0.04662833454583212"
$ws.Range("F22").Value = [double]"-0.003895637640006719"
$ws.Range("G22").Value = [double]"-0.001039033487937697"
$ws.Range("H22").Value = [double]"0.09562662117888582"
$ws.Range("B23").Value = [double]"-0.00503613637582422"
$ws.Range("C23").Value = [double]"0.0005112007405269867"
$ws.Range("D23").Value = [double]"-9.390414256769018"
$ws.Range("E23").Value = [double]"1.019289903783663e-08"
$ws.Range("F23").Value = [double]"-0.006038074584400107"
$ws.Range("G23").Value = [double]"-0.004034198167248333"
$ws.Range("H23").Value = [double]"0.09305782036703381"
$ws.Range("B24").Value = [double]"-0.00415510939198079"
$ws.Range("C24").Value = [double]"0.0005077619555907203"
$ws.Range("D24").Value = [double]"-7.322245901547114"
$ws.Range("E24").Value = [double]"0.01834297707943118"
$ws.Range("F24").Value = [double]"-0.005150307678227482"
$ws.Range("G24").Value = [double]"-0.003159911105734098"
$ws.Range("H24").Value = [double]"0.09393884735087724"
$ws.Range("B25").Value = [double]"-0.002965001684836056"
$ws.Range("C25").Value = [double]"0.0004984298957980975"
$ws.Range("D25").Value = [double]"-5.474408928284702"
$ws.Range("E25").Value = [double]"0.1234674631563146"
$ws.Range("F25").Value = [double]"-0.003941909419209885"
$ws.Range("G25").Value = [double]"-0.001988093950462227"
$ws.Range("H25").Value = [double]"0.09512895505802196"
$ws.Range("B26").Value = [double]"0.01689858836606458"
$ws.Range("C26").Value = [double]"0.001559237756845363"
$ws.Range("D26").Value = [double]"20.10854211153268"
$ws.Range("E26").Value = [double]"0.02323030177506833"
$ws.Range("F26").Value = [double]"0.01384252851717928"
$ws.Range("G26").Value = [double]"0.01995464821494989"
$ws.Range("H26").Value = [double]"0.1149925451089226"

Write-Output "Updated 168 cells"
